$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift data up by one row (drop the oldest record) and recompute
# Trad_Prediction (C) as the previous Real_Close (B), with fresh
# AI_Prediction (D) simulation values, for rows 2-35.
$ws.Cells.Item(2, 1).Value = 45847.60416666666
$ws.Cells.Item(2, 2).Value = 501.5299987792969
$ws.Cells.Item(2, 3).Value = 504.5599975585938
$ws.Cells.Item(2, 4).Value = 509.3401283152386
$ws.Cells.Item(3, 1).Value = 45847.64583333334
$ws.Cells.Item(3, 2).Value = 501.0899963378906
$ws.Cells.Item(3, 3).Value = 501.5299987792969
$ws.Cells.Item(3, 4).Value = 501.2669041153314
$ws.Cells.Item(4, 1).Value = 45847.6875
$ws.Cells.Item(4, 2).Value = 501.4700012207031
$ws.Cells.Item(4, 3).Value = 501.0899963378906
$ws.Cells.Item(4, 4).Value = 485.4786069838178
$ws.Cells.Item(5, 1).Value = 45847.72916666666
$ws.Cells.Item(5, 2).Value = 502.7950134277344
$ws.Cells.Item(5, 3).Value = 501.4700012207031
$ws.Cells.Item(5, 4).Value = 490.8370997962443
$ws.Cells.Item(6, 1).Value = 45847.77083333334
$ws.Cells.Item(6, 2).Value = 502.3200073242188
$ws.Cells.Item(6, 3).Value = 502.7950134277344
$ws.Cells.Item(6, 4).Value = 507.9550117837441
$ws.Cells.Item(7, 1).Value = 45847.8125
$ws.Cells.Item(7, 2).Value = 503.4800109863281
$ws.Cells.Item(7, 3).Value = 502.3200073242188
$ws.Cells.Item(7, 4).Value = 508.1956426674561
$ws.Cells.Item(8, 1).Value = 45848.5625
$ws.Cells.Item(8, 2).Value = 498.5400085449219
$ws.Cells.Item(8, 3).Value = 503.4800109863281
$ws.Cells.Item(8, 4).Value = 505.0552318309559
$ws.Cells.Item(9, 1).Value = 45848.60416666666
$ws.Cells.Item(9, 2).Value = 500.0650024414062
$ws.Cells.Item(9, 3).Value = 498.5400085449219
$ws.Cells.Item(9, 4).Value = 514.9930552887106
$ws.Cells.Item(10, 1).Value = 45848.64583333334
$ws.Cells.Item(10, 2).Value = 500.2749938964844
$ws.Cells.Item(10, 3).Value = 500.0650024414062
$ws.Cells.Item(10, 4).Value = 515.8148026627798
$ws.Cells.Item(11, 1).Value = 45848.6875
$ws.Cells.Item(11, 2).Value = 500.9150085449219
$ws.Cells.Item(11, 3).Value = 500.2749938964844
$ws.Cells.Item(11, 4).Value = 487.5638192880039
$ws.Cells.Item(12, 1).Value = 45848.72916666666
$ws.Cells.Item(12, 2).Value = 501.3900146484375
$ws.Cells.Item(12, 3).Value = 500.9150085449219
$ws.Cells.Item(12, 4).Value = 517.9709313178091
$ws.Cells.Item(13, 1).Value = 45848.77083333334
$ws.Cells.Item(13, 2).Value = 501.5499877929688
$ws.Cells.Item(13, 3).Value = 501.3900146484375
$ws.Cells.Item(13, 4).Value = 497.8013294701834
$ws.Cells.Item(14, 1).Value = 45848.8125
$ws.Cells.Item(14, 2).Value = 501.5199890136719
$ws.Cells.Item(14, 3).Value = 501.5499877929688
$ws.Cells.Item(14, 4).Value = 504.1304629740764
$ws.Cells.Item(15, 1).Value = 45849.5625
$ws.Cells.Item(15, 2).Value = 501.5
$ws.Cells.Item(15, 3).Value = 501.5199890136719
$ws.Cells.Item(15, 4).Value = 506.606505034066
$ws.Cells.Item(16, 1).Value = 45849.60416666666
$ws.Cells.Item(16, 2).Value = 502.7900085449219
$ws.Cells.Item(16, 3).Value = 501.5
$ws.Cells.Item(16, 4).Value = 498.8253337325588
$ws.Cells.Item(17, 1).Value = 45849.64583333334
$ws.Cells.Item(17, 2).Value = 504.7099914550781
$ws.Cells.Item(17, 3).Value = 502.7900085449219
$ws.Cells.Item(17, 4).Value = 484.770510397758
$ws.Cells.Item(18, 1).Value = 45849.6875
$ws.Cells.Item(18, 2).Value = 504.5450134277344
$ws.Cells.Item(18, 3).Value = 504.7099914550781
$ws.Cells.Item(18, 4).Value = 498.821212957596
$ws.Cells.Item(19, 1).Value = 45849.72916666666
$ws.Cells.Item(19, 2).Value = 504.2799987792969
$ws.Cells.Item(19, 3).Value = 504.5450134277344
$ws.Cells.Item(19, 4).Value = 511.142941753954
$ws.Cells.Item(20, 1).Value = 45849.77083333334
$ws.Cells.Item(20, 2).Value = 503.4028930664062
$ws.Cells.Item(20, 3).Value = 504.2799987792969
$ws.Cells.Item(20, 4).Value = 507.0127189260152
$ws.Cells.Item(21, 1).Value = 45849.8125
$ws.Cells.Item(21, 2).Value = 503.1700134277344
$ws.Cells.Item(21, 3).Value = 503.4028930664062
$ws.Cells.Item(21, 4).Value = 500.3138399580718
$ws.Cells.Item(22, 1).Value = 45852.5625
$ws.Cells.Item(22, 2).Value = 501.5499877929688
$ws.Cells.Item(22, 3).Value = 503.1700134277344
$ws.Cells.Item(22, 4).Value = 487.6974701485398
$ws.Cells.Item(23, 1).Value = 45852.60416666666
$ws.Cells.Item(23, 2).Value = 502.9649963378906
$ws.Cells.Item(23, 3).Value = 501.5499877929688
$ws.Cells.Item(23, 4).Value = 491.515116754712
$ws.Cells.Item(24, 1).Value = 45852.64583333334
$ws.Cells.Item(24, 2).Value = 502.4800109863281
$ws.Cells.Item(24, 3).Value = 502.9649963378906
$ws.Cells.Item(24, 4).Value = 508.1424275296509
$ws.Cells.Item(25, 1).Value = 45852.6875
$ws.Cells.Item(25, 2).Value = 502.7200012207031
$ws.Cells.Item(25, 3).Value = 502.4800109863281
$ws.Cells.Item(25, 4).Value = 507.6075053080272
$ws.Cells.Item(26, 1).Value = 45852.72916666666
$ws.Cells.Item(26, 2).Value = 502.8999938964844
$ws.Cells.Item(26, 3).Value = 502.7200012207031
$ws.Cells.Item(26, 4).Value = 516.8012665661424
$ws.Cells.Item(27, 1).Value = 45852.77083333334
$ws.Cells.Item(27, 2).Value = 502.9700012207031
$ws.Cells.Item(27, 3).Value = 502.8999938964844
$ws.Cells.Item(27, 4).Value = 500.26636444873
$ws.Cells.Item(28, 1).Value = 45852.8125
$ws.Cells.Item(28, 2).Value = 502.9400024414062
$ws.Cells.Item(28, 3).Value = 502.9700012207031
$ws.Cells.Item(28, 4).Value = 512.6847099647705
$ws.Cells.Item(29, 1).Value = 45853.5625
$ws.Cells.Item(29, 2).Value = 505.6499938964844
$ws.Cells.Item(29, 3).Value = 502.9400024414062
$ws.Cells.Item(29, 4).Value = 510.0913777821144
$ws.Cells.Item(30, 1).Value = 45853.60416666666
$ws.Cells.Item(30, 2).Value = 507.2799987792969
$ws.Cells.Item(30, 3).Value = 505.6499938964844
$ws.Cells.Item(30, 4).Value = 504.6393796146909
$ws.Cells.Item(31, 1).Value = 45853.64583333334
$ws.Cells.Item(31, 2).Value = 507.3200073242188
$ws.Cells.Item(31, 3).Value = 507.2799987792969
$ws.Cells.Item(31, 4).Value = 500.4710672097632
$ws.Cells.Item(32, 1).Value = 45853.6875
$ws.Cells.Item(32, 2).Value = 506.8699951171875
$ws.Cells.Item(32, 3).Value = 507.3200073242188
$ws.Cells.Item(32, 4).Value = 511.9772331720771
$ws.Cells.Item(33, 1).Value = 45853.72916666666
$ws.Cells.Item(33, 2).Value = 507.9500122070312
$ws.Cells.Item(33, 3).Value = 506.8699951171875
$ws.Cells.Item(33, 4).Value = 510.2332390427385
$ws.Cells.Item(34, 1).Value = 45853.77083333334
$ws.Cells.Item(34, 2).Value = 506.5799865722656
$ws.Cells.Item(34, 3).Value = 507.9500122070312
$ws.Cells.Item(34, 4).Value = 484.489156672649
$ws.Cells.Item(35, 1).Value = 45853.8125
$ws.Cells.Item(35, 2).Value = 505.6300048828125
$ws.Cells.Item(35, 3).Value = 506.5799865722656
$ws.Cells.Item(35, 4).Value = 506.6558908424671

# Remove the now-obsolete last row (was row 36, data fully shifted up)
$ws.Rows(36).Delete()
